# Apply the worksheet edits described by the commit diff:
#   - C2: 11   -> 12
#   - C3: 9.5  -> 11
#   - C4: 1.4  -> 1.45
#   - selected cell moves from B3 to C3
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("C2").Value = 12
$ws.Range("C3").Value = 11
$ws.Range("C4").Value = 1.45

# Reflect the updated active selection (was B3, now C3)
[void]$ws.Range("C3").Select()
